$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The "Pina/Caramelo" Vega Monumental Concepcion price-history block
# (rows 135-208) shifts up by one row: row N now holds what used to be
# in row N+1. A brand-new observation is inserted at the top of the block
# (row 135) and the record that used to be the very last one (old row 208)
# is now duplicated down into a freshly appended row 209. Columns
# A,B,C,E,F,G,H,I,J are identical for every data row in this block, so
# only D (col 4) and K:T (cols 11-20) change.
# -----------------------------------------------------------------------
$data = @{}
$data[135] = @(44839, "Caramelo", "Segunda", 230, 22000, 23000, 22652, "`$/caja 14 unidades", "Ecuador", 1618, 14)
$data[136] = @(44327, "Caramelo", "Segunda", 200, 16500, 17000, 16750, "`$/caja 14 unidades", "Ecuador", 1196, 14)
$data[137] = @(44460, "Caramelo", "Segunda", 150, 19000, 20000, 19667, "`$/caja 14 unidades", "Ecuador", 1405, 14)
$data[138] = @(44546, "Caramelo", "Primera", 100, 17000, 18000, 17500, "`$/caja 12 unidades", "Ecuador", 1458, 12)
$data[139] = @(44609, "Caramelo", "Segunda", 200, 16000, 17000, 16500, "`$/caja 14 unidades", "Ecuador", 1179, 14)
$data[140] = @(44308, "Caramelo", "Segunda", 200, 16000, 16500, 16250, "`$/caja 14 unidades", "Ecuador", 1161, 14)
$data[141] = @(44776, "Caramelo", "Primera", 200, 19000, 20000, 19500, "`$/caja 12 unidades", "Ecuador", 1625, 12)
$data[142] = @(44218, "Caramelo", "Primera", 200, 14000, 15000, 14500, "`$/caja 12 unidades", "Ecuador", 1208, 12)
$data[143] = @(44680, "Caramelo", "Primera", 220, 16000, 17000, 16545, "`$/caja 14 unidades", "Ecuador", 1182, 14)
$data[144] = @(44498, "Caramelo", "Segunda", 200, 19000, 20000, 19500, "`$/caja 14 unidades", "Ecuador", 1393, 14)
$data[145] = @(44159, "Caramelo", "Segunda", 200, 21000, 22000, 21500, "`$/caja 14 unidades", "Ecuador", 1536, 14)
$data[146] = @(44215, "Caramelo", "Segunda", 400, 15000, 16000, 15500, "`$/caja 14 unidades", "Ecuador", 1107, 14)
$data[147] = @(44210, "Caramelo", "Primera", 200, 14000, 15000, 14500, "`$/caja 12 unidades", "Ecuador", 1208, 12)
$data[148] = @(44238, "Caramelo", "Primera", 400, 15500, 16000, 15750, "`$/caja 12 unidades", "Ecuador", 1312, 12)
$data[149] = @(44749, "Caramelo", "Segunda", 200, 17000, 18000, 17500, "`$/caja 10 unidades", "Ecuador", 1750, 10)
$data[150] = @(44224, "Caramelo", "Primera", 100, 16000, 17000, 16500, "`$/caja 12 unidades", "Ecuador", 1375, 12)
$data[151] = @(44811, "Caramelo", "Segunda", 350, 17000, 18000, 17429, "`$/caja 14 unidades", "Ecuador", 1245, 14)
$data[152] = @(44166, "Caramelo", "Segunda", 200, 22000, 23000, 22500, "`$/caja 14 unidades", "Ecuador", 1607, 14)
$data[153] = @(44567, "Caramelo", "Segunda", 310, 14000, 15000, 14484, "`$/caja 14 unidades", "Ecuador", 1035, 14)
$data[154] = @(44614, "Caramelo", "Segunda", 200, 16000, 17000, 16500, "`$/caja 14 unidades", "Ecuador", 1179, 14)
$data[155] = @(44435, "Caramelo", "Primera", 100, 18000, 18500, 18250, "`$/caja 12 unidades", "Ecuador", 1521, 12)
$data[156] = @(44435, "Caramelo", "Segunda", 400, 18000, 18500, 18250, "`$/caja 14 unidades", "Ecuador", 1304, 14)
$data[157] = @(44264, "Caramelo", "Segunda", 200, 15500, 16000, 15750, "`$/caja 14 unidades", "Ecuador", 1125, 14)
$data[158] = @(44442, "Caramelo", "Segunda", 200, 19000, 19500, 19250, "`$/caja 14 unidades", "Ecuador", 1375, 14)
$data[159] = @(44755, "Sin especificar", "Segunda", 270, 18000, 19000, 18556, "`$/caja 14 unidades", "Ecuador", 1325, 14)
$data[160] = @(44336, "Caramelo", "Segunda", 200, 16000, 16500, 16250, "`$/caja 14 unidades", "Ecuador", 1161, 14)
$data[161] = @(44426, "Caramelo", "Segunda", 200, 18000, 19000, 18500, "`$/caja 14 unidades", "Ecuador", 1321, 14)
$data[162] = @(44826, "Sin especificar", "Tercera", 220, 20000, 21000, 20545, "`$/caja 16 unidades", "Ecuador", 1284, 16)
$data[163] = @(44670, "Caramelo", "Segunda", 200, 15000, 16000, 15500, "`$/caja 14 unidades", "Ecuador", 1107, 14)
$data[164] = @(44161, "Caramelo", "Primera", 100, 20000, 21000, 20500, "`$/caja 12 unidades", "Ecuador", 1708, 12)
$data[165] = @(44161, "Caramelo", "Segunda", 200, 20000, 21000, 20500, "`$/caja 14 unidades", "Ecuador", 1464, 14)
$data[166] = @(44448, "Caramelo", "Primera", 200, 18500, 19000, 18750, "`$/caja 12 unidades", "Ecuador", 1562, 12)
$data[167] = @(44253, "Caramelo", "Segunda", 400, 15500, 16000, 15750, "`$/caja 14 unidades", "Ecuador", 1125, 14)
$data[168] = @(44789, "Sin especificar", "Segunda", 200, 18000, 19000, 18500, "`$/caja 14 unidades", "Ecuador", 1321, 14)
$data[169] = @(44484, "Caramelo", "Primera", 100, 20000, 21000, 20500, "`$/caja 12 unidades", "Ecuador", 1708, 12)
$data[170] = @(44645, "Caramelo", "Segunda", 180, 16000, 17000, 16444, "`$/caja 14 unidades", "Ecuador", 1175, 14)
$data[171] = @(44595, "Caramelo", "Segunda", 180, 14000, 15000, 14444, "`$/caja 14 unidades", "Ecuador", 1032, 14)
$data[172] = @(44630, "Caramelo", "Primera", 200, 19000, 20000, 19500, "`$/caja 12 unidades", "Ecuador", 1625, 12)
$data[173] = @(44784, "Caramelo", "Segunda", 220, 18000, 19000, 18545, "`$/caja 14 unidades", "Ecuador", 1325, 14)
$data[174] = @(44679, "Caramelo", "Segunda", 180, 15000, 16000, 15556, "`$/caja 14 unidades", "Ecuador", 1111, 14)
$data[175] = @(44635, "Caramelo", "Segunda", 110, 16000, 17000, 16545, "`$/caja 14 unidades", "Bolivia", 1182, 14)
$data[176] = @(44334, "Caramelo", "Segunda", 100, 16000, 16500, 16250, "`$/caja 14 unidades", "Ecuador", 1161, 14)
$data[177] = @(44334, "Caramelo", "Tercera", 100, 16000, 16500, 16250, "`$/caja 16 unidades", "Ecuador", 1016, 16)
$data[178] = @(44722, "Caramelo", "Primera", 270, 16000, 17000, 16556, "`$/caja 14 unidades", "Ecuador", 1183, 14)
$data[179] = @(44706, "Caramelo", "Segunda", 200, 16000, 17000, 16500, "`$/caja 14 unidades", "Ecuador", 1179, 14)
$data[180] = @(44610, "Caramelo", "Segunda", 150, 15000, 16000, 15533, "`$/caja 14 unidades", "Ecuador", 1110, 14)
$data[181] = @(44196, "Caramelo", "Segunda", 200, 15000, 16000, 15500, "`$/caja 14 unidades", "Ecuador", 1107, 14)
$data[182] = @(44194, "Caramelo", "Primera", 200, 15000, 16000, 15500, "`$/caja 12 unidades", "Ecuador", 1292, 12)
$data[183] = @(44446, "Caramelo", "Segunda", 200, 18000, 19000, 18500, "`$/caja 14 unidades", "Ecuador", 1321, 14)
$data[184] = @(44273, "Caramelo", "Segunda", 200, 15000, 15500, 15250, "`$/caja 14 unidades", "Ecuador", 1089, 14)
$data[185] = @(44672, "Caramelo", "Primera", 200, 15000, 16000, 15500, "`$/caja 12 unidades", "Ecuador", 1292, 12)
$data[186] = @(44274, "Caramelo", "Primera", 150, 15500, 16000, 15833, "`$/caja 12 unidades", "Ecuador", 1319, 12)
$data[187] = @(44782, "Sin especificar", "Segunda", 270, 18000, 19000, 18444, "`$/caja 14 unidades", "Ecuador", 1317, 14)
$data[188] = @(44421, "Caramelo", "Primera", 100, 18000, 18000, 18000, "`$/caja 12 unidades", "Ecuador", 1500, 12)
$data[189] = @(44421, "Caramelo", "Segunda", 200, 17000, 18000, 17500, "`$/caja 14 unidades", "Ecuador", 1250, 14)
$data[190] = @(44376, "Caramelo", "Segunda", 200, 16500, 17000, 16750, "`$/caja 14 unidades", "Ecuador", 1196, 14)
$data[191] = @(44390, "Caramelo", "Segunda", 200, 15000, 16000, 15500, "`$/caja 14 unidades", "Ecuador", 1107, 14)
$data[192] = @(44726, "Caramelo", "Segunda", 200, 17000, 18000, 17500, "`$/caja 14 unidades", "Ecuador", 1250, 14)
$data[193] = @(44357, "Caramelo", "Primera", 200, 16000, 16500, 16250, "`$/caja 12 unidades", "Ecuador", 1354, 12)
$data[194] = @(44747, "Caramelo", "Primera", 100, 19000, 20000, 19500, "`$/caja 14 unidades", "Ecuador", 1393, 14)
$data[195] = @(44244, "Caramelo", "Segunda", 200, 14500, 15000, 14750, "`$/caja 14 unidades", "Ecuador", 1054, 14)
$data[196] = @(44433, "Caramelo", "Primera", 100, 18000, 18500, 18250, "`$/caja 12 unidades", "Ecuador", 1521, 12)
$data[197] = @(44385, "Caramelo", "Segunda", 200, 17000, 17500, 17250, "`$/caja 14 unidades", "Ecuador", 1232, 14)
$data[198] = @(44162, "Caramelo", "Segunda", 200, 20000, 21000, 20500, "`$/caja 14 unidades", "Ecuador", 1464, 14)
$data[199] = @(44204, "Caramelo", "Primera", 200, 14500, 15000, 14750, "`$/caja 12 unidades", "Ecuador", 1229, 12)
$data[200] = @(44355, "Caramelo", "Segunda", 200, 16000, 16500, 16250, "`$/caja 14 unidades", "Ecuador", 1161, 14)
$data[201] = @(44763, "Caramelo", "Segunda", 200, 18000, 19000, 18500, "`$/caja 14 unidades", "Ecuador", 1321, 14)
$data[202] = @(44565, "Caramelo", "Primera", 250, 15000, 16000, 15400, "`$/caja 14 unidades", "Ecuador", 1100, 14)
$data[203] = @(44187, "Caramelo", "Segunda", 200, 15000, 16000, 15500, "`$/caja 14 unidades", "Ecuador", 1107, 14)
$data[204] = @(44775, "Caramelo", "Segunda", 200, 18000, 19000, 18500, "`$/caja 14 unidades", "Bolivia", 1321, 14)
$data[205] = @(44575, "Caramelo", "Segunda", 200, 16000, 16500, 16250, "`$/caja 14 unidades", "Ecuador", 1161, 14)
$data[206] = @(44553, "Caramelo", "Primera", 220, 14000, 15000, 14455, "`$/caja 14 unidades", "Ecuador", 1032, 14)
$data[207] = @(44292, "Caramelo", "Segunda", 400, 15500, 16000, 15750, "`$/caja 14 unidades", "Ecuador", 1125, 14)
$data[208] = @(44453, "Caramelo", "Segunda", 200, 21000, 22000, 21500, "`$/caja 14 unidades", "Ecuador", 1536, 14)
$data[209] = @(44832, "Caramelo", "Segunda", 150, 21000, 22000, 21667, "`$/caja 14 unidades", "Ecuador", 1548, 14)

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    # Column D (Fecha)
    $ws.Cells.Item($row, 4).Value2 = $vals[0]

    # Columns K:T (Variedad .. Kg/unidad) -- contiguous block, cols 11-20
    $arr = New-Object "object[,]" 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $vals[$i + 1] }
    $ws.Range($ws.Cells.Item($row, 11), $ws.Cells.Item($row, 20)).Value2 = $arr
}

# Row 209 is brand new -- it needs the constant identifying columns too
# (A,B,C,E,F,G,H,I,J), which are the same on every row of this block.
$ws.Range("A209").Value2 = 11
$ws.Range("B209").Value2 = "Vega Monumental Concepción"
$ws.Range("C209").Value2 = "Bíobío"
$ws.Range("E209").Value2 = 8
$ws.Range("F209").Value2 = "Fruta"
$ws.Range("G209").Value2 = 100108
$ws.Range("H209").Value2 = "Tropicales y subtropicales"
$ws.Range("I209").Value2 = 100108005
$ws.Range("J209").Value2 = "Piña"

# Match the date-formatted style already used by the rest of column D.
$ws.Range("D209").NumberFormat = $ws.Range("D208").NumberFormat

